$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price / volume data scraped on 2023-02-07 (GitHub Actions refresh).
# Column D (Price) and E (Volume 1h) are stored as plain text in the sheet, so force
# the "@" text number format before assigning to avoid Excel coercing them to numbers.
$data = @(
    @{ Cell = "D2"; Value = "328.81" }
    @{ Cell = "E2"; Value = "0.24%" }
    @{ Cell = "D3"; Value = "44.22" }
    @{ Cell = "E3"; Value = "1.03%" }
    @{ Cell = "D4"; Value = "5.484" }
    @{ Cell = "E4"; Value = "-0.99%" }
    @{ Cell = "D5"; Value = "0.08046" }
    @{ Cell = "E5"; Value = "0.35%" }
    @{ Cell = "D6"; Value = "2.029" }
    @{ Cell = "E6"; Value = "7.05%" }
    @{ Cell = "D7"; Value = "0.9539" }
    @{ Cell = "E7"; Value = "1.00%" }
    @{ Cell = "D8"; Value = "0.1111" }
    @{ Cell = "E8"; Value = "-6.52%" }
    @{ Cell = "D9"; Value = "0.1876" }
    @{ Cell = "E10"; Value = "1.11%" }
    @{ Cell = "D11"; Value = "0.09986" }
    @{ Cell = "E11"; Value = "3.06%" }
    @{ Cell = "D12"; Value = "0.04734" }
    @{ Cell = "E12"; Value = "5.68%" }
    @{ Cell = "D13"; Value = "0.1060" }
    @{ Cell = "E13"; Value = "-0.73%" }
    @{ Cell = "D14"; Value = "0.001256" }
    @{ Cell = "E14"; Value = "-1.42%" }
    @{ Cell = "D15"; Value = "0.04086" }
    @{ Cell = "E15"; Value = "-2.58%" }
    @{ Cell = "D16"; Value = "0.005767" }
    @{ Cell = "E16"; Value = "-2.85%" }
    @{ Cell = "D17"; Value = "3.376" }
    @{ Cell = "E17"; Value = "-0.89%" }
    @{ Cell = "D18"; Value = "4.417" }
    @{ Cell = "E18"; Value = "3.74%" }
    @{ Cell = "D19"; Value = "2.658" }
    @{ Cell = "E19"; Value = "3.83%" }
    @{ Cell = "E20"; Value = "-0.61%" }
    @{ Cell = "D21"; Value = "0.1400" }
    @{ Cell = "E21"; Value = "-0.99%" }
    @{ Cell = "E22"; Value = "2.94%" }
    @{ Cell = "E23"; Value = "5.06%" }
    @{ Cell = "D24"; Value = "0.004340" }
    @{ Cell = "E24"; Value = "0.73%" }
    @{ Cell = "E25"; Value = "-0.78%" }
    @{ Cell = "E26"; Value = "-6.23%" }
    @{ Cell = "D38"; Value = "0.02578" }
    @{ Cell = "E38"; Value = "-2.80%" }
    @{ Cell = "D39"; Value = "0.05667" }
    @{ Cell = "E39"; Value = "3.45%" }
    @{ Cell = "D40"; Value = "0.007739" }
    @{ Cell = "E40"; Value = "2.07%" }
    @{ Cell = "D41"; Value = "0.1398" }
    @{ Cell = "E41"; Value = "0.25%" }
    @{ Cell = "D42"; Value = "0.007362" }
    @{ Cell = "E42"; Value = "-9.52%" }
    @{ Cell = "D43"; Value = "0.002011" }
    @{ Cell = "E43"; Value = "0.38%" }
    @{ Cell = "D44"; Value = "0.008520" }
    @{ Cell = "E44"; Value = "-3.26%" }
    @{ Cell = "D45"; Value = "0.00007086" }
    @{ Cell = "E45"; Value = "-0.42%" }
    @{ Cell = "E46"; Value = "0.01%" }
    @{ Cell = "D47"; Value = "0.0005807" }
    @{ Cell = "E47"; Value = "-0.08%" }
    @{ Cell = "B48"; Value = "CoinbaseStockToken" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin" }
    @{ Cell = "D48"; Value = "0.003503" }
    @{ Cell = "E48"; Value = "54.12%" }
    @{ Cell = "B49"; Value = "BOLO" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo" }
    @{ Cell = "D49"; Value = "0.003513" }
    @{ Cell = "E49"; Value = "3.52%" }
    @{ Cell = "E50"; Value = "0.01%" }
    @{ Cell = "E51"; Value = "0.01%" }
)

foreach ($item in $data) {
    $cell = $ws.Range($item.Cell)
    if ($item.Cell[0] -eq "D" -or $item.Cell[0] -eq "E") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $item.Value
}
